$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-07-13 20:57:46"

for ($r = 2; $r -le 73; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}
